$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 160.71428
$ws.Range("I33").Value = 121.40909
$ws.Range("J33").Value = 304.83334
$ws.Range("K33").Value = 121.40909
$ws.Range("L33").Value = 304.83334
$ws.Range("M33").Value = 107.59091
$ws.Range("N33").Value = -762.83334
$ws.Range("H116").Value = 3442.8235
$ws.Range("I116").Value = 2684
$ws.Range("K116").Value = 2684
$ws.Range("M116").Value = 758
$ws.Range("H135").Value = 244.71428
$ws.Range("I135").Value = 216.85715
$ws.Range("J135").Value = 272.57144
$ws.Range("K135").Value = 1951.71435
$ws.Range("L135").Value = 2453.14296
$ws.Range("M135").Value = 583.28565
$ws.Range("N135").Value = -7523.14296
$ws.Range("H136").Value = 69000
$ws.Range("J136").Value = 69000
$ws.Range("L136").Value = 69000
$ws.Range("N136").Value = -79200

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3539.8696
$ws.Range("J32").Value = 1894.6666
$ws.Range("L32").Value = 1894.6666
$ws.Range("N32").Value = -2468.6666
$ws.Range("H61").Value = 1699.1818
$ws.Range("I61").Value = 1334.625
$ws.Range("K61").Value = 1334.625
$ws.Range("M61").Value = -1122.625
$ws.Range("H136").Value = 1699.1818
$ws.Range("I136").Value = 1334.625
$ws.Range("K136").Value = 4003.875
$ws.Range("M136").Value = -1453.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5090.28
$ws.Range("I134").Value = 1052
$ws.Range("J134").Value = 34704.332
$ws.Range("K134").Value = 3156
$ws.Range("L134").Value = 104112.996
$ws.Range("M134").Value = -621
$ws.Range("N134").Value = -109182.996

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 277.8889
$ws.Range("I7").Value = 302.625
$ws.Range("K7").Value = 302.625
$ws.Range("M7").Value = -189.625
$ws.Range("H31").Value = 1089.7084
$ws.Range("I31").Value = 748.61536
$ws.Range("J31").Value = 1976.55
$ws.Range("K31").Value = 748.61536
$ws.Range("L31").Value = 1976.55
$ws.Range("M31").Value = -453.61536
$ws.Range("N31").Value = -2566.55
$ws.Range("H34").Value = 1089.7084
$ws.Range("I34").Value = 748.61536
$ws.Range("J34").Value = 1976.55
$ws.Range("K34").Value = 748.61536
$ws.Range("L34").Value = 1976.55
$ws.Range("M34").Value = -546.61536
$ws.Range("N34").Value = -2380.55
$ws.Range("H58").Value = 696.5
$ws.Range("I58").Value = 696.5
$ws.Range("K58").Value = 696.5
$ws.Range("M58").Value = -493.5
$ws.Range("H94").Value = 850.5714
$ws.Range("I94").Value = 1870.3334
$ws.Range("J94").Value = 572.4545000000001
$ws.Range("K94").Value = 1870.3334
$ws.Range("L94").Value = 572.4545000000001
$ws.Range("M94").Value = -1419.3334
$ws.Range("N94").Value = -1474.4545
$ws.Range("H105").Value = 670
$ws.Range("I105").Value = 670
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 670
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 1077
$ws.Range("N105").ClearContents()
$ws.Range("H107").Value = 608.95
$ws.Range("I107").Value = 159.71428
$ws.Range("J107").Value = 850.8461
$ws.Range("K107").Value = 159.71428
$ws.Range("L107").Value = 850.8461
$ws.Range("M107").Value = 1760.28572
$ws.Range("N107").Value = -4690.8461
$ws.Range("H134").Value = 6803914
$ws.Range("I134").Value = 8772996
$ws.Range("K134").Value = 26318988
$ws.Range("M134").Value = -26316453
$ws.Range("H136").Value = 696.5
$ws.Range("I136").Value = 696.5
$ws.Range("K136").Value = 2089.5
$ws.Range("M136").Value = 460.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 2780
$ws.Range("I22").Value = 2600
$ws.Range("J22").Value = 3500
$ws.Range("K22").Value = 7800
$ws.Range("L22").Value = 10500
$ws.Range("M22").Value = -7631
$ws.Range("N22").Value = -10838
$ws.Range("H27").Value = 2780
$ws.Range("I27").Value = 2600
$ws.Range("J27").Value = 3500
$ws.Range("K27").Value = 7800
$ws.Range("L27").Value = 10500
$ws.Range("M27").Value = -7698
$ws.Range("N27").Value = -10704
$ws.Range("H33").Value = 250.21428
$ws.Range("I33").Value = 200.33333
$ws.Range("J33").Value = 340
$ws.Range("K33").Value = 1201.99998
$ws.Range("L33").Value = 2040
$ws.Range("M33").Value = -918.9999800000001
$ws.Range("N33").Value = -2606
$ws.Range("H38").Value = 230
$ws.Range("J38").Value = 400
$ws.Range("L38").Value = 1200
$ws.Range("N38").Value = -1894
$ws.Range("H39").Value = 2000
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 2000
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 6000
$ws.Range("N39").Value = -6588
$ws.Range("M39").ClearContents()
$ws.Range("H49").Value = 1000
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 1000
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 3000
$ws.Range("N49").Value = -3312
$ws.Range("M49").ClearContents()
$ws.Range("H92").Value = 661.6
$ws.Range("J92").Value = 661.6
$ws.Range("L92").Value = 1984.8
$ws.Range("N92").Value = -4480.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3041.4
$ws.Range("I126").Value = 1916.4445
$ws.Range("J126").Value = 3961.818
$ws.Range("K126").Value = 5749.333500000001
$ws.Range("L126").Value = 11885.454
$ws.Range("M126").Value = -3279.333500000001
$ws.Range("N126").Value = -16825.454
$ws.Range("H132").Value = 2770.3333
$ws.Range("I132").Value = 2531.7368
$ws.Range("K132").Value = 7595.2104
$ws.Range("M132").Value = -5065.2104
$ws.Range("H135").Value = 49998.668
$ws.Range("J135").Value = 49998.668
$ws.Range("L135").Value = 49998.668
$ws.Range("N135").Value = -60138.668

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 10288.5
$ws.Range("I136").Value = 14023.25
$ws.Range("K136").Value = 42069.75
$ws.Range("M136").Value = -39519.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 606.63635
$ws.Range("I107").Value = 607.7778
$ws.Range("K107").Value = 1823.3334
$ws.Range("M107").Value = 96.66660000000002
$ws.Range("H136").Value = 1357.6666
$ws.Range("I136").Value = 504
$ws.Range("K136").Value = 1512
$ws.Range("M136").Value = 1038
